$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add the new "PHOTO" column header in F1 (extends the used range to A1:F11)
$ws.Range("F1").Value = "PHOTO"

# The rest of the new column (F2:F11) holds blank/empty text values, one per
# existing data row. Plain `.Value = ""` collapses a cell to fully blank
# (Excel drops an empty-string write), so each cell is entered as the
# literal formula ="" which evaluates to, and is typed as, an empty string -
# the closest achievable equivalent of an empty text cell for every row.
for ($row = 2; $row -le 11; $row++) {
    $ws.Cells.Item($row, 6).Formula = "="""""
}
